$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "eyJraWQiOiJVS3ZSSnBseHZOQUJLVDhhWkZtMzZ3ZnpXenB3cHV6dmloczNob0tWbDlvPSIsImFsZyI6IlJTMjU2In0.eyJzdWIiOiJiYmE0OGQ0NC01N2Q3LTQzZGYtYmY1Yi04MmY5YmE4YzYyODciLCJpc3MiOiJodHRwczpcL1wvY29nbml0by1pZHAudXMtZWFzdC0xLmFtYXpvbmF3cy5jb21cL3VzLWVhc3QtMV82SUljcUM0WFAiLCJjbGllbnRfaWQiOiIyaW9zb2U4OTUzY2owaWhtbzdydWt2bWg4Iiwib3JpZ2luX2p0aSI6IjZiMDlkYjA1LTk1NjgtNGNjZS1hYTZiLTM3ZDQ5ZGZjY2VkOCIsImV2ZW50X2lkIjoiYTkyZDMyMjktNTJhYS00NzI0LWE1N2UtMzRjOTU4OGY4MTc5IiwidG9rZW5fdXNlIjoiYWNjZXNzIiwic2NvcGUiOiJhd3MuY29nbml0by5zaWduaW4udXNlci5hZG1pbiIsImF1dGhfdGltZSI6MTY1ODkyNDMxNSwiZXhwIjoxNjU4OTI1MjE1LCJpYXQiOjE2NTg5MjQzMTUsImp0aSI6IjNlMWYxODU2LTRjYWUtNDIzZi1iMTUzLTZjNTg2MzNlYWI3NiIsInVzZXJuYW1lIjoiYmJhNDhkNDQtNTdkNy00M2RmLWJmNWItODJmOWJhOGM2Mjg3In0.dlgudeNMDbcctJe4qwn53KDADkxBvD2wp4K5wfIlWA-HY2So2UCc9wzIq6jnt9pCFNW2GU_559zvm_mO8fef9dUEeNr0wGwLvew9k1q8eOlBhZAu2A4931-JvGyzNNeWkeoGK9rmHHnY2-6rIkUmxoJJVfCWVxEX3EzvIrJi9wUtL4TP5FeeGX_YNqLy60etC_zijGKvdQVMm5RDCQR7m8UZMly8cOAIMEstxnMFU99MYm5ibY2WGWQA4sJe3EFYsl4oU_2Csc684U5YJkR8m90R5qrBGUX2OMvKUgWxa5fp-HZHV-YlvlssbmmW26-OF8SjPWOKimpDELiknou4dA"
$ws.Range("C3").Value = "eyJraWQiOiIzUWJLOXJkWTBlbUdcLyt0VjZHbUVVaEVsSXZTb2JqdFVDcTFNQnZ5SDFnST0iLCJhbGciOiJSUzI1NiJ9.eyJzdWIiOiJiYmE0OGQ0NC01N2Q3LTQzZGYtYmY1Yi04MmY5YmE4YzYyODciLCJlbWFpbF92ZXJpZmllZCI6dHJ1ZSwiaXNzIjoiaHR0cHM6XC9cL2NvZ25pdG8taWRwLnVzLWVhc3QtMS5hbWF6b25hd3MuY29tXC91cy1lYXN0LTFfNklJY3FDNFhQIiwiY3VzdG9tOnVzZXJfaWQiOiJ1c3JfbDJva3VhaDlVb3ZtUTJZNiIsImNvZ25pdG86dXNlcm5hbWUiOiJiYmE0OGQ0NC01N2Q3LTQzZGYtYmY1Yi04MmY5YmE4YzYyODciLCJjdXN0b206YWNjb3VudHMiOiJ7XCJhY2N0X2wyb2t3MjVvNUVkcjM2ckFcIjpcIk93bmVyXCIsXCJhY2N0X2wydTdkN3M4VXFEYmo3SzJcIjpcIk93bmVyXCIsXCJhY2N0X2wzbjA1NGtmME5Pb3lPb25cIjpcIk93bmVyXCIsXCJhY2N0X2wydTdmaDhwOEI3VjVZdXVcIjpcIk93bmVyXCIsXCJhY2N0X2wydTdnNGh6a29vVTR4WHBcIjpcIk93bmVyXCIsXCJhY2N0X2wzbjAyb200OXRhRkZuYWtcIjpcIk93bmVyXCIsXCJhY2N0X2wzbjA0YmF5OHBsOHVLcXBcIjpcIk93bmVyXCIsXCJhY2N0X2wzbjA3bm9rSE9RQUdiWEJcIjpcIk93bmVyXCIsXCJhY2N0X2wydTdmNWgzUlpxdFdqTE1cIjpcIk93bmVyXCIsXCJhY2N0X2wzbjA0NXMwS1kzTzB5eHdcIjpcIk93bmVyXCIsXCJhY2N0X2wyb2t1YWhic3F2M2N6dHpcIjpcIk93bmVyXCIsXCJhY2N0X2wydTdmMXplRklEUkxiaVVcIjpcIk93bmVyXCIsXCJhY2N0X2wzbjA3cGU5SU04eUs4d2JcIjpcIk93bmVyXCIsXCJhY2N0X2wzbXl0bHlqaml1R21ueDFcIjpcIk93bmVyXCIsXCJhY2N0X2wzbjA3a28zcWNRa25qN3pcIjpcIk93bmVyXCIsXCJhY2N0X2wyc3Nob2JsTmt4cXhPMUdcIjpcIk93bmVyXCIsXCJhY2N0X2wzbjA3MzFoSVUxWm1tTTVcIjpcIk93bmVyXCIsXCJhY2N0X2w0M3FlOGJiTjNjN0N5NmlcIjpcIk93bmVyXCIsXCJhY2N0X2wydTZkY2F0clhGTlJ6R1dcIjpcIk93bmVyXCIsXCJhY2N0X2wzbjA1MDBnNFhZak42bWFcIjpcIk93bmVyXCIsXCJhY2N0X2wydTdnNnNuc1lPdFRwSmJcIjpcIk93bmVyXCJ9Iiwib3JpZ2luX2p0aSI6IjZiMDlkYjA1LTk1NjgtNGNjZS1hYTZiLTM3ZDQ5ZGZjY2VkOCIsImF1ZCI6IjJpb3NvZTg5NTNjajBpaG1vN3J1a3ZtaDgiLCJldmVudF9pZCI6ImE5MmQzMjI5LTUyYWEtNDcyNC1hNTdlLTM0Yzk1ODhmODE3OSIsInRva2VuX3VzZSI6ImlkIiwiYXV0aF90aW1lIjoxNjU4OTI0MzE1LCJuYW1lIjoiUGFsYXNoIFBhdGlkYXIiLCJleHAiOjE2NTg5MjUyMTUsImlhdCI6MTY1ODkyNDMxNSwianRpIjoiMTk5MmM0YmUtYTEyYS00MjI2LTlhMTAtNjRkMmJjNWUwYmYxIiwiZW1haWwiOiJwYWxhc2hAdHJ5c3BlZWQuY29tIn0.DOQGcF_K0XvAae2vtLgxYxO9mfXPNAYmnG2S-bZpghtDVMB4A4LzbNwl9DaVdN7PWoarlGQirVTWDOi1q2PiHEoLHDi4OoOrd7meALyrqFJOGGeBv0B3GSv3IDA6-txxX5T53Vr_8FwjDLB-PIiKK57_PEaoI1_BDiQPThl3SJ6KCRl3NNC7sw-bJpJwPdPwk9LuMuGlduQeAbsV-BLofBR7HVSdN8Wi3daxuzkd9h3oVBAASU1HoaUrH0dk2AVRVL3pI0JYCPMfEDxA7oaBHUVRhaPr3PydLFtPP9rYPf1USp6CCD0_ZkFC9oDqQRtItr7Wyk2q06nKeCb2YVbujg"
$ws.Range("D3").Value = "eyJjdHkiOiJKV1QiLCJlbmMiOiJBMjU2R0NNIiwiYWxnIjoiUlNBLU9BRVAifQ.ENnm_jL6sUy0Hqio9pIXePXMnXK3Czkne7-S98-4wLtApbu_1jaRKb3mpDDmpK8QRc0lr0uwy9NY_ttPijGXyUwmiHKdpp3ZRC0qOI4D_dbtcX9I0XjgsczBjuFguhIblELPUt0OjLUcPL2lA2DmC96xKUzxOaeFfUixtKecqtD9Y-FpBTKtjlGuLELk-DQBLMXVCar4ltJDQ04wpp1WQqQZi0k6FkjsjR0F4CbWj2iDQ33NMiV2Zl5ikRYB-BAdeHNb5QG7xiog0re222NB9CMpPZx5uYl1_Q5fk5JvwSZz38rI9dyovm1j56vQPtCCksbaXtah5WGY-sY4UfNWxA.Npy0vCj6Em7v_jKB.-LTYd86xPSFa4D7qmrew9eNgHUEC9IXw6UMYzV1SuM43UQmpZ4VlMbQsj3kBpW0PEC0ZFuT2xi-LHxjHVmTsplQZ1e1cqg-myNBd9l78R2Q57jo0F695h-ZoSOuPwhEuvujZT2vPsEb7WzswQM56xsaPfHqkAstq1BBeB1ubVhbfl28f6_lE46ijipw0KUAGEOVeGqBEMir9rFMU19mFkiOR02PFzaEhJBpLF_1uxVjfiyC_I2yGE9pQ0josgxMigHnoxnicShf96ZkTTi8O0xipwchEEq17FSkW_76ILV_WdHbVtIguHc4oKf_xmdT-zDP650iH-crSYkxephQv9qLMt6Qu6WOCkttS1462y-a46GfyMzcyPtvK5Nnlptnxnlv63lBszU1tXepd6SnkGENYm3F7i61QrjbTF81t1B4wAtiJtDqo6tk-4aqtHNOkmxqbutJ1OM_lnyDybTWGlY1xIAyOoIsR5-4mt8dedSk3ZB6mR4gIwENTEonKILbwAo312rQKKJUw4ckMFXzVjHS63sxIYe-WuOtK0UrbmgkEgkwUyCorruv89F6H0EQxyqpe_KGPTON0o-XLR79ONoQjwQ9B0EmTKiGp6gQpqADpkggtQpy3CbbD0ATFHq5tyqDHfpOotlIofHbcllfE5E0V1_tcsDvkVfeq_PsDsH-GwIG9emc25iGNltIek_cUgzptFRuEY2C-5_7E5OC6jTJR8amMNRcg6fzoAQBznClkPFlmxoeNpNzLdao563wEbQZ3GYxK_HGnbpS6Obxn_IN6r2BogoNf_MI_QrBdl4YsIWYalzAC7rX6xod3LGGW3B-4AJxLcaudY81qhgcFlNlRtkVIZiCFfm4e1R4eqx11ih7oARu9IBrOZjrVNRktUH7Pqc6NvBYqmYkkgjruDjIQ_pOqchNg8t6Epi-Ytzb4iW1F7ZgifnE-axx2EeHyhxUgW7BJTdecIab659BfCXR9vLQDFUzgvCmTM9YNTpjSjKNAEoyKERhF8bnMVKWtMW6t6NEQLjxi7NVNdxVoqjD51JAFSQ2UL1T2u1WA71sY1AjUnihFssbv_qWmNEVkb3jyvz8s6zyysRQy8JSiIQ0Ezk3sCPoyL-GoSpn_8aHxIn3p3vxfgLtWSg1hxx4pObrkncq2ate20H16m3J1ziF7vJtm5yW7oamACWXPOfCn0N3dNIn3uOrJxDg3_ANj3ca2UQSel22ggt6xt6kRsC3iaRRbBhp1RcdKmuuJf0PmkT2VGxLBjmaboK9ZHCVfxh0GQTFdPJjQIlR92szeMke_aTSBIxmdIa-YbdZQEaWdpDCFBmgk7uhc.klbxqBh_gWfCvIZX_Lp9Dg"
$ws.Range("E3").Value = "sn_l63kpdwlViAzSY1q"
